# Revert "Bcrypt integration for secure password hashing - Part 2"
#
# This undoes the status changes that Part 2 made to the task tracker and
# restores the "Transition backend..." row to Done / clears the bcrypt
# row's status / reinstates the "Add 'Profile' page" row as Done, then
# re-adds a dropped reference row pointing at the hashing article.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add "Profile" page" was "In progress" -> now "Done"
$ws.Range("D5").Value = "Done"

# "Use bcrypt for hashing and storing passwords*" status cleared (was "Done")
$ws.Range("D14").ClearContents()

# "Transition backend to use database instead of static memory" was
# "In progress" -> now "Done"
$ws.Range("D15").Value = "Done"

# Re-add the reference link row that was previously removed
$ws.Range("A31").Value = "https://crackstation.net/hashing-security.htm"

# Restore the view/selection state left behind by the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D21").Select()
